$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2290287.2
$ws.Range("J40").Value = 3625638.2
$ws.Range("L40").Value = 3625638.2
$ws.Range("N40").Value = -3625988.2
$ws.Range("H136").Value = 22367.691
$ws.Range("J136").Value = 22367.691
$ws.Range("L136").Value = 22367.691
$ws.Range("N136").Value = -32567.691
$ws.Range("H137").Value = 25038.35
$ws.Range("I137").Value = 1598.7858
$ws.Range("J137").Value = 68792.2
$ws.Range("K137").Value = 4796.357400000001
$ws.Range("L137").Value = 206376.6
$ws.Range("M137").Value = -2246.357400000001
$ws.Range("N137").Value = -211476.6
$ws.Range("H138").Value = 10871886
$ws.Range("J138").Value = 2552.4365
$ws.Range("L138").Value = 7657.309499999999
$ws.Range("N138").Value = -17937.3095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1817.1765
$ws.Range("I2").Value = 1680.0869
$ws.Range("J2").Value = 2103.818
$ws.Range("K2").Value = 1680.0869
$ws.Range("L2").Value = 2103.818
$ws.Range("M2").Value = -1567.0869
$ws.Range("N2").Value = -2329.818
$ws.Range("H45").Value = 3774.4443
$ws.Range("I45").Value = 3594.3076
$ws.Range("J45").Value = 4242.8
$ws.Range("K45").Value = 3594.3076
$ws.Range("L45").Value = 4242.8
$ws.Range("M45").Value = -3217.3076
$ws.Range("N45").Value = -4996.8
$ws.Range("H110").Value = 481.64706
$ws.Range("I110").Value = 497
$ws.Range("J110").Value = 444.8
$ws.Range("K110").Value = 497
$ws.Range("L110").Value = 444.8
$ws.Range("M110").Value = 1548
$ws.Range("N110").Value = -4534.8
$ws.Range("H116").Value = 1817.1765
$ws.Range("I116").Value = 1680.0869
$ws.Range("J116").Value = 2103.818
$ws.Range("K116").Value = 1680.0869
$ws.Range("L116").Value = 2103.818
$ws.Range("M116").Value = 613.9131
$ws.Range("N116").Value = -6691.818
$ws.Range("H122").Value = 1459.3529
$ws.Range("I122").Value = 1158.8975
$ws.Range("J122").Value = 2435.8333
$ws.Range("K122").Value = 3476.6925
$ws.Range("L122").Value = 7307.499899999999
$ws.Range("M122").Value = -1026.6925
$ws.Range("N122").Value = -12207.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1817.1765
$ws.Range("I3").Value = 1680.0869
$ws.Range("J3").Value = 2103.818
$ws.Range("K3").Value = 1680.0869
$ws.Range("L3").Value = 2103.818
$ws.Range("M3").Value = -1566.0869
$ws.Range("N3").Value = -2331.818
$ws.Range("H22").Value = 1000458.2
$ws.Range("I22").Value = 1250347.8
$ws.Range("K22").Value = 1250347.8
$ws.Range("M22").Value = -1250174.8
$ws.Range("H94").Value = 900
$ws.Range("I94").Value = 900
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -449
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 2027.7778
$ws.Range("I105").Value = 1545.4546
$ws.Range("J105").Value = 2785.7144
$ws.Range("K105").Value = 1545.4546
$ws.Range("L105").Value = 2785.7144
$ws.Range("M105").Value = 201.5454
$ws.Range("N105").Value = -6279.7144
$ws.Range("H107").Value = 1366.0416
$ws.Range("I107").Value = 1290.9333
$ws.Range("J107").Value = 1491.2222
$ws.Range("K107").Value = 1290.9333
$ws.Range("L107").Value = 1491.2222
$ws.Range("M107").Value = 629.0667000000001
$ws.Range("N107").Value = -5331.2222
$ws.Range("H134").Value = 66318.12
$ws.Range("I134").Value = 66318.12
$ws.Range("K134").Value = 198954.36
$ws.Range("M134").Value = -196419.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12254.091
$ws.Range("I31").Value = 15610.682
$ws.Range("J31").Value = 5540.909
$ws.Range("K31").Value = 15610.682
$ws.Range("L31").Value = 5540.909
$ws.Range("M31").Value = -15315.682
$ws.Range("N31").Value = -6130.909
$ws.Range("H34").Value = 12254.091
$ws.Range("I34").Value = 15610.682
$ws.Range("J34").Value = 5540.909
$ws.Range("K34").Value = 15610.682
$ws.Range("L34").Value = 5540.909
$ws.Range("M34").Value = -15408.682
$ws.Range("N34").Value = -5944.909
$ws.Range("H99").Value = 22731846
$ws.Range("I99").Value = 4256.1875
$ws.Range("J99").Value = 83338750
$ws.Range("K99").Value = 4256.1875
$ws.Range("L99").Value = 83338750
$ws.Range("M99").Value = -2758.1875
$ws.Range("N99").Value = -83341746
$ws.Range("H126").Value = 22731846
$ws.Range("I126").Value = 4256.1875
$ws.Range("J126").Value = 83338750
$ws.Range("K126").Value = 12768.5625
$ws.Range("L126").Value = 250016250
$ws.Range("M126").Value = -10298.5625
$ws.Range("N126").Value = -250021190

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3546.4285
$ws.Range("I114").Value = 2481.5
$ws.Range("J114").Value = 4966.3335
$ws.Range("K114").Value = 7444.5
$ws.Range("L114").Value = 14899.0005
$ws.Range("M114").Value = -4190.5
$ws.Range("N114").Value = -21407.0005
$ws.Range("H118").Value = 38468570
$ws.Range("I118").Value = 100000216
$ws.Range("J118").Value = 11287.125
$ws.Range("K118").Value = 300000648
$ws.Range("L118").Value = 33861.375
$ws.Range("M118").Value = -299999405
$ws.Range("N118").Value = -36347.375
$ws.Range("H119").Value = 3241.5557
$ws.Range("I119").Value = 1938.1666
$ws.Range("K119").Value = 5814.4998
$ws.Range("M119").Value = -976.4997999999996
$ws.Range("H122").Value = 1234.6923
$ws.Range("I122").Value = 349.14285
$ws.Range("J122").Value = 1560.9474
$ws.Range("K122").Value = 3142.28565
$ws.Range("L122").Value = 14048.5266
$ws.Range("M122").Value = -692.2856500000003
$ws.Range("N122").Value = -18948.5266
$ws.Range("H138").Value = 1213.3334
$ws.Range("I138").Value = 1213.3334
$ws.Range("K138").Value = 3640.0002
$ws.Range("M138").Value = 1499.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13167.5
$ws.Range("I70").Value = 10738.462
$ws.Range("J70").Value = 17678.572
$ws.Range("K70").Value = 10738.462
$ws.Range("L70").Value = 17678.572
$ws.Range("M70").Value = -10468.462
$ws.Range("N70").Value = -18218.572
$ws.Range("H73").Value = 13167.5
$ws.Range("I73").Value = 10738.462
$ws.Range("J73").Value = 17678.572
$ws.Range("K73").Value = 10738.462
$ws.Range("L73").Value = 17678.572
$ws.Range("M73").Value = -9802.462
$ws.Range("N73").Value = -19550.572
$ws.Range("H102").Value = 20834714
$ws.Range("I102").Value = 23810802
$ws.Range("J102").Value = 2104.6667
$ws.Range("K102").Value = 23810802
$ws.Range("L102").Value = 2104.6667
$ws.Range("M102").Value = -23809180
$ws.Range("N102").Value = -5348.6667
$ws.Range("H113").Value = 1318.3871
$ws.Range("I113").Value = 1044.6666
$ws.Range("J113").Value = 2256.8572
$ws.Range("K113").Value = 1044.6666
$ws.Range("L113").Value = 2256.8572
$ws.Range("M113").Value = 1125.3334
$ws.Range("N113").Value = -6596.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 16000
$ws.Range("I45").Value = 13000
$ws.Range("J45").Value = 17000
$ws.Range("K45").Value = 13000
$ws.Range("L45").Value = 17000
$ws.Range("M45").Value = -12593
$ws.Range("N45").Value = -17814
$ws.Range("H46").Value = 746.125
$ws.Range("I46").Value = 737.3333
$ws.Range("J46").Value = 757.4286
$ws.Range("K46").Value = 737.3333
$ws.Range("L46").Value = 757.4286
$ws.Range("M46").Value = -549.3333
$ws.Range("N46").Value = -1133.4286
